$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.344.82'
$ws.Range("E2").Value = '  -3.54%  '

# Row 3
$ws.Range("D3").Value = '2.248.64'
$ws.Range("E3").Value = '  -4.25%  '

# Row 4
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.53'
$ws.Range("E5").Value = '  -2.89%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.634'
$ws.Range("E6").Value = '  -5.45%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.81'
$ws.Range("E7").Value = '  -3.43%  '

# Row 8
$ws.Range("E8").Value = '  +0.14%  '

# Row 9
$ws.Range("E9").Value = '  -5.69%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0994'
$ws.Range("E10").Value = '  -1.02%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.32'
$ws.Range("E11").Value = '  +0.04%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.75'
$ws.Range("E12").Value = '  +8.68%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  -3.91%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.82'
$ws.Range("E14").Value = '  -5.94%  '

# Row 15
$ws.Range("D15").Value = '2.581.58'

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.08'
$ws.Range("E16").Value = '  -7.46%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.862'
$ws.Range("E17").Value = '  -4.49%  '

# Row 18
$ws.Range("D18").Value = '2.249.38'
$ws.Range("E18").Value = '  -4.33%  '

# Row 19
$ws.Range("D19").Value = '42.179.24'
$ws.Range("E19").Value = '  -3.72%  '

# Row 20
$ws.Range("E20").Value = '  -4.72%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.26'
$ws.Range("E21").Value = '  -6.80%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.50'
$ws.Range("E22").Value = '  -6.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.13'
$ws.Range("E23").Value = '  -6.87%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.02'
$ws.Range("E24").Value = '  +4.12%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.66'
$ws.Range("E26").Value = '  -2.24%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.36'
$ws.Range("E27").Value = '  -5.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.07'
$ws.Range("E28").Value = '  -3.82%  '

# Row 29
$ws.Range("E29").Value = '  -4.08%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '169.31'
$ws.Range("E30").Value = '  -4.34%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.67'
$ws.Range("E31").Value = '  -7.96%  '

# Row 32
$ws.Range("E32").Value = '  -6.21%  '

# Row 33
$ws.Range("E33").Value = '  -6.66%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.37'
$ws.Range("E34").Value = '  -0.69%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0716'
$ws.Range("E35").Value = '  -4.28%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.77'
$ws.Range("E36").Value = '  -7.00%  '

# Row 37
$ws.Range("E37").Value = '  -3.59%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.11'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.26'
$ws.Range("E39").Value = '  -5.03%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.06'
$ws.Range("E40").Value = '  -5.71%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0266'
$ws.Range("E41").Value = '  -3.25%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '66.71'
$ws.Range("E42").Value = '  -1.55%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.98'
$ws.Range("E43").Value = '  -3.45%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.02'
$ws.Range("E44").Value = '  -1.96%  '

# Row 45
$ws.Range("E45").Value = '  -3.69%  '

# Row 46
$ws.Range("E46").Value = '  -7.80%  '

# Row 47
$ws.Range("E47").Value = '  +0.13%  '

# Row 48
$ws.Range("B48").Value = 'SynthetixNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.39'
$ws.Range("E48").Value = '  +7.49%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.37'
$ws.Range("E49").Value = '  -4.94%  '

# Row 50
$ws.Range("B50").Value = 'TrustWalletToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.19'
$ws.Range("E50").Value = '  -4.44%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.14'
$ws.Range("E51").Value = '  +7.39%  '
